# Update "想去人数" (F column) figures on the 展览 / 演出 / 全部类型 sheets
# to the refreshed counts captured at commit 456a3b4 (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 17
$ws.Range("F3").Value = 20969
$ws.Range("F4").Value = 815
$ws.Range("F6").Value = 1120
$ws.Range("F7").Value = 24
$ws.Range("F8").Value = 7829
$ws.Range("F9").Value = 547
$ws.Range("F11").Value = 754
$ws.Range("F12").Value = 299
$ws.Range("F14").Value = 183
$ws.Range("F15").Value = 156
$ws.Range("F16").Value = 29
$ws.Range("F18").Value = 217
$ws.Range("F19").Value = 1354
$ws.Range("F20").Value = 504
$ws.Range("F24").Value = 80
$ws.Range("F25").Value = 82
$ws.Range("F26").Value = 344
$ws.Range("F27").Value = 1168
$ws.Range("F29").Value = 35
$ws.Range("F30").Value = 211
$ws.Range("F31").Value = 5215
$ws.Range("F32").Value = 596
$ws.Range("F33").Value = 127
$ws.Range("F34").Value = 4986
$ws.Range("F37").Value = 30
$ws.Range("F39").Value = 12973
$ws.Range("F40").Value = 1357
$ws.Range("F41").Value = 118
$ws.Range("F42").Value = 46
$ws.Range("F43").Value = 67
$ws.Range("F45").Value = 412
$ws.Range("F46").Value = 4046
$ws.Range("F47").Value = 327

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 322

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 17
$ws.Range("F3").Value = 20969
$ws.Range("F4").Value = 815
$ws.Range("F6").Value = 1120
$ws.Range("F7").Value = 24
$ws.Range("F8").Value = 7829
$ws.Range("F9").Value = 547
$ws.Range("F11").Value = 754
$ws.Range("F12").Value = 299
$ws.Range("F14").Value = 183
$ws.Range("F15").Value = 156
$ws.Range("F16").Value = 29
$ws.Range("F17").Value = 217
$ws.Range("F18").Value = 1354
$ws.Range("F19").Value = 504
$ws.Range("F23").Value = 80
$ws.Range("F24").Value = 82
$ws.Range("F25").Value = 344
$ws.Range("F26").Value = 1168
$ws.Range("F28").Value = 35
$ws.Range("F29").Value = 211
$ws.Range("F30").Value = 322
$ws.Range("F31").Value = 596
$ws.Range("F33").Value = 127
$ws.Range("F35").Value = 4986
$ws.Range("F38").Value = 30
$ws.Range("F40").Value = 12973
$ws.Range("F41").Value = 1357
$ws.Range("F42").Value = 118
$ws.Range("F43").Value = 46
$ws.Range("F44").Value = 67
$ws.Range("F46").Value = 412
$ws.Range("F47").Value = 4046
$ws.Range("F48").Value = 327
